# Implements: PowerPoint template and slide duplication logic
#
# The deck has 10 "slots" (slide1..slide10), each a fixed template
# (one caption TextBox + 4 cropped-picture shapes) showing one source
# image's info ("sample_image_NN.jpeg"). Two more source images needed
# slots, so the template gets duplicated twice (from the last two
# slides) and appended at the end, and every slot's caption/crop
# descriptions shift down by two, with the first two slots going blank.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Step 1: duplicate the last two template slides (9 and 10) so we have
# two fresh slides to hold the overflow content, and place them at the
# very end of the deck (positions 11 and 12).
# ---------------------------------------------------------------------

# Duplicate slide 9 -> lands right after it (position 10), pushing the
# original slide 10 to position 11. Move the new duplicate to the end.
$dupOf9 = $p.Slides.Item(9).Duplicate()
$p.Slides.Item(10).MoveTo($p.Slides.Count)

# Original slide 10 is back at position 10. Duplicate it -> lands at
# position 11, pushing the slide-9 duplicate to position 12.
$dupOf10 = $p.Slides.Item(10).Duplicate()

# Now order is: ... 9(orig), 10(orig), 11(dup of 10), 12(dup of 9).
# Swap the last two so the dup-of-9 ("09") precedes dup-of-10 ("10").
$p.Slides.Item(11).MoveTo($p.Slides.Count)

# Final base order: 1..8 unchanged, 9=orig9, 10=orig10,
# 11=dup-of-9, 12=dup-of-10.

# ---------------------------------------------------------------------
# Step 2: helper to relabel a slide's caption + the descr (alt text) on
# each of its 4 picture shapes.
# ---------------------------------------------------------------------
function Set-SlideContent($slideIndex, $caption, $descrs) {
    $slide = $p.Slides.Item($slideIndex)
    $captionShape = $slide.Shapes.Item(1)
    $tr = $captionShape.TextFrame.TextRange
    # Replace the run text via Characters() (rather than TextRange.Text
    # directly) so PowerPoint doesn't stamp an explicit <a:rPr lang="..."/>
    # on the run; then restore the caption box height since the
    # spAutoFit body recomputes it to the new text's natural size.
    $tr.Characters(1, $tr.Text.Length).Text = $caption
    $captionShape.Height = 36
    for ($i = 0; $i -lt 4; $i++) {
        $slide.Shapes.Item(2 + $i).AlternativeText = $descrs[$i]
    }
}

Set-SlideContent 3  "sample_image_01.jpeg" @("tmpz7qawo4c.png","tmpwb89m5iq.png","tmpvymg3nw8.png","tmptzlgphh2.png")
Set-SlideContent 4  "sample_image_02.jpeg" @("tmpdqx3zl0a.png","tmpq18pnc_2.png","tmp0an413fe.png","tmp79f76vgb.png")
Set-SlideContent 5  "sample_image_03.jpeg" @("tmpjvelttvd.png","tmpu8mnjok9.png","tmp3uvnn_yv.png","tmpyrl1geyr.png")
Set-SlideContent 6  "sample_image_04.jpeg" @("tmp229ygck_.png","tmp5u32cay8.png","tmpbliuao7w.png","tmpi2caycy_.png")
Set-SlideContent 7  "sample_image_05.jpeg" @("tmp_abigt05.png","tmp5lnqdqow.png","tmpw1yrzwtw.png","tmpo_qkj1r0.png")
Set-SlideContent 8  "sample_image_06.jpeg" @("tmpfd040epp.png","tmpzxa_kkzb.png","tmp7yhfo_dd.png","tmppk7tvzj2.png")
Set-SlideContent 9  "sample_image_07.jpeg" @("tmp1mrer_dg.png","tmpgupedown.png","tmpczl607xc.png","tmpomb7nlgz.png")
Set-SlideContent 10 "sample_image_08.jpeg" @("tmpsq30y_cj.png","tmpq99ggyxk.png","tmp2z8jhr5d.png","tmppv8g9aie.png")
Set-SlideContent 11 "sample_image_09.jpeg" @("tmpo0wte7v8.png","tmpfwxy203u.png","tmp4wbvhxwp.png","tmppmeo4m0_.png")
Set-SlideContent 12 "sample_image_10.jpeg" @("tmpjav419d9.png","tmpnk3o4gmh.png","tmpku74kuta.png","tmpxpfn37dk.png")

# ---------------------------------------------------------------------
# Step 3: the first two template slots no longer have content to show
# (their former images shifted down), so clear them to blank slides.
# ---------------------------------------------------------------------
foreach ($idx in 1, 2) {
    $slide = $p.Slides.Item($idx)
    while ($slide.Shapes.Count -gt 0) {
        $slide.Shapes.Item(1).Delete()
    }
}
